# Applies the "Clean up of templates, incorporating dukes 5.1 into standard
# processing methods" edit described by the diff:
#   * fix a handful of mislabeled unit / item / category cells on sheet 7.2
#     (adds a new shared string "Electrical outout" along the way)
#   * re-autofit every sheet's columns (now including column A) and drop the
#     zoom level each sheet was left at back down to a sane value
#   * leave the workbook focused on sheet 7.2 instead of 7.7

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Data corrections on "7.2" (Total CHP schemes listing).
#    Row 31 had unit "number" where it should have read "MW".
#    Rows 60-88 had a stray "Electrical capacity" / "Electrical output
#    (GWh) [note 3]" / "MW" combination copy-pasted down the column instead
#    of the correct "Electrical outout" / "All" / "GWh" triple.
# ---------------------------------------------------------------------------
$wsDukes = $wb.Worksheets.Item("7.2")

$wsDukes.Range("E31").Value = "MW"

for ($r = 60; $r -le 88; $r++) {
    $wsDukes.Range("C$r").Value = "Electrical outout"
    $wsDukes.Range("D$r").Value = "All"
    $wsDukes.Range("E$r").Value = "GWh"
}

# ---------------------------------------------------------------------------
# 2. Per-sheet view clean-up: autofit every column (including the row-number
#    column A, which previously relied on the default width) and settle on
#    the new zoom level. Activating each sheet also resets its selection
#    back to A1 and clears any stored topLeftCell scroll position.
# ---------------------------------------------------------------------------
$zooms = [ordered]@{
    "7.1.A" = 100
    "7.1.B" = 115
    "7.3.A" = 100
    "7.3.B" = 115
    "7.3.C" = 115
    "7.3.D" = 100
    "7.3.E" = 115
    "7.3.F" = 100
    "7.4.A" = 115
    "7.4.B" = 130
    "7.7"   = 115
    "7.8.A" = 130
    "7.8.B" = 100
}

foreach ($name in $zooms.Keys) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()
    $excel.ActiveWindow.Zoom = $zooms[$name]
    $ws.UsedRange.EntireColumn.AutoFit()
    $ws.Range("A1").Select()
}

# "7.2" is autofitted/zoomed/activated last so it ends up the workbook's
# active tab (activeTab moves from 7.7 to 7.2), matching the saved view.
$wsDukes.Activate()
$excel.ActiveWindow.Zoom = 115
$wsDukes.UsedRange.EntireColumn.AutoFit()
$wsDukes.Range("A1").Select()
